$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 2419.6
$ws.Range("I38").Value = 699.6667
$ws.Range("J38").Value = 4999.5
$ws.Range("K38").Value = 2099.0001
$ws.Range("L38").Value = 14998.5
$ws.Range("M38").Value = -1727.0001
$ws.Range("N38").Value = -15742.5
# Row 40
$ws.Range("H40").Value = 254291.5
$ws.Range("I40").Value = 2687.375
$ws.Range("K40").Value = 2687.375
$ws.Range("M40").Value = -2512.375
# Row 62
$ws.Range("H62").Value = 9987
$ws.Range("I62").Value = 9986.5
$ws.Range("K62").Value = 9986.5
$ws.Range("M62").Value = -9362.5
# Row 65
$ws.Range("H65").Value = 9987
$ws.Range("I65").Value = 9986.5
$ws.Range("K65").Value = 49932.5
$ws.Range("M65").Value = -46812.5
# Row 74
$ws.Range("H74").Value = 10423.286
$ws.Range("I74").Value = 10423.286
$ws.Range("K74").Value = 10423.286
$ws.Range("M74").Value = -9487.286
# Row 76
$ws.Range("H76").Value = 3780.4
$ws.Range("I76").Value = 3634.3333
$ws.Range("J76").Value = 3999.5
$ws.Range("K76").Value = 3634.3333
$ws.Range("L76").Value = 3999.5
$ws.Range("M76").Value = -3319.3333
$ws.Range("N76").Value = -4629.5
# Row 77
$ws.Range("H77").Value = 10423.286
$ws.Range("I77").Value = 10423.286
$ws.Range("K77").Value = 52116.43
$ws.Range("M77").Value = -47436.43
# Row 79
$ws.Range("H79").Value = 3780.4
$ws.Range("I79").Value = 3634.3333
$ws.Range("J79").Value = 3999.5
$ws.Range("K79").Value = 3634.3333
$ws.Range("L79").Value = 3999.5
$ws.Range("M79").Value = -2542.3333
$ws.Range("N79").Value = -6183.5
# Row 86
$ws.Range("H86").Value = 4062.2727
$ws.Range("J86").Value = 3835.625
$ws.Range("L86").Value = 3835.625
$ws.Range("N86").Value = -6081.625
# Row 89
$ws.Range("H89").Value = 4062.2727
$ws.Range("J89").Value = 3835.625
$ws.Range("L89").Value = 19178.125
$ws.Range("N89").Value = -30410.125
# Row 98
$ws.Range("H98").Value = 4142.727
$ws.Range("I98").Value = 3083.1428
$ws.Range("K98").Value = 3083.1428
$ws.Range("M98").Value = -1585.1428
# Row 112
$ws.Range("H112").Value = 2577.9092
$ws.Range("J112").Value = 2777.7778
$ws.Range("L112").Value = 8333.3334
$ws.Range("N112").Value = -10549.3334
# Row 122
$ws.Range("H122").Value = 4142.727
$ws.Range("I122").Value = 3083.1428
$ws.Range("K122").Value = 9249.428400000001
$ws.Range("M122").Value = -6799.428400000001
# Row 138
$ws.Range("H138").Value = 4646.325
$ws.Range("I138").Value = 2636.125
$ws.Range("J138").Value = 5148.875
$ws.Range("K138").Value = 7908.375
$ws.Range("L138").Value = 15446.625
$ws.Range("M138").Value = -2768.375
$ws.Range("N138").Value = -25726.625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8491.058999999999
$ws.Range("I32").Value = 8095.1875
$ws.Range("J32").Value = 14825
$ws.Range("K32").Value = 8095.1875
$ws.Range("L32").Value = 14825
$ws.Range("M32").Value = -7808.1875
$ws.Range("N32").Value = -15399
# Row 61
$ws.Range("H61").Value = 2275.6667
$ws.Range("I61").Value = 2244.1765
$ws.Range("K61").Value = 2244.1765
$ws.Range("M61").Value = -2032.1765
# Row 136
$ws.Range("H136").Value = 2275.6667
$ws.Range("I136").Value = 2244.1765
$ws.Range("K136").Value = 6732.529500000001
$ws.Range("M136").Value = -4182.529500000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 130
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 175
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 175
$ws.Range("M22").Value = 73
$ws.Range("N22").Value = -521
# Row 134
$ws.Range("H134").Value = 1478.7
$ws.Range("I134").Value = 1458.25
$ws.Range("K134").Value = 4374.75
$ws.Range("M134").Value = -1839.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 665.7143
$ws.Range("I22").Value = 536.6
$ws.Range("K22").Value = 536.6
$ws.Range("M22").Value = -186.6
# Row 31
$ws.Range("H31").Value = 2043.4546
$ws.Range("I31").Value = 1164.2222
$ws.Range("J31").Value = 6000
$ws.Range("K31").Value = 1164.2222
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = -869.2221999999999
$ws.Range("N31").Value = -6590
# Row 34
$ws.Range("H34").Value = 2043.4546
$ws.Range("I34").Value = 1164.2222
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 1164.2222
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -962.2221999999999
$ws.Range("N34").Value = -6404
# Row 58
$ws.Range("H58").Value = 3701.1
$ws.Range("I58").Value = 2258.7144
$ws.Range("K58").Value = 2258.7144
$ws.Range("M58").Value = -2055.7144
# Row 136
$ws.Range("H136").Value = 3701.1
$ws.Range("I136").Value = 2258.7144
$ws.Range("K136").Value = 6776.1432
$ws.Range("M136").Value = -4226.1432

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 188.42857
$ws.Range("I2").Value = 192
$ws.Range("K2").Value = 192
$ws.Range("M2").Value = -79
# Row 122
$ws.Range("H122").Value = 3140.4285
$ws.Range("I122").Value = 2997.1667
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8991.500100000001
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -6541.500100000001
$ws.Range("N122").Value = -16900

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 8021.3
$ws.Range("I22").Value = 6471
$ws.Range("K22").Value = 6471
$ws.Range("M22").Value = -6176
# Row 27
$ws.Range("H27").Value = 8021.3
$ws.Range("I27").Value = 6471
$ws.Range("K27").Value = 6471
$ws.Range("M27").Value = -6364
# Row 46
$ws.Range("H46").Value = 1666.6666
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -2376
# Row 55
$ws.Range("H55").Value = 1259.2
$ws.Range("I55").Value = 1247.5
$ws.Range("J55").Value = 1262.125
$ws.Range("K55").Value = 1247.5
$ws.Range("L55").Value = 1262.125
$ws.Range("M55").Value = -1074.5
$ws.Range("N55").Value = -1608.125
# Row 61
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -798
$ws.Range("N61").ClearContents()
# Row 100
$ws.Range("H100").Value = 996.6
$ws.Range("I100").Value = 996.6
$ws.Range("K100").Value = 996.6
$ws.Range("M100").Value = -455.6
# Row 113
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1170
$ws.Range("N113").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 39130
$ws.Range("J75").Value = 39130
$ws.Range("L75").Value = 39130
$ws.Range("N75").Value = -41002
# Row 78
$ws.Range("H78").Value = 39130
$ws.Range("J78").Value = 39130
$ws.Range("L78").Value = 117390
$ws.Range("N78").Value = -126750
# Row 113
$ws.Range("H113").Value = 650
$ws.Range("J113").Value = 620.5
$ws.Range("L113").Value = 1861.5
$ws.Range("N113").Value = -6201.5
